$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.583.35'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '3.334.80'
$ws.Range('E3').Value = '  -0.04%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '586.85'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '183.14'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.69%  '
$ws.Range('E7').Value = '  +4.09%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '3.338.84'
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('E10').Value = '  -0.92%  '
$ws.Range('E11').Value = '  +2.35%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.404'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.51%  '
$ws.Range('D13').Value = '3.918.89'
$ws.Range('E13').Value = '  +0.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.131'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.56%  '
$ws.Range('D15').Value = '66.580.45'
$ws.Range('E15').Value = '  +0.01%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '26.63'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.31%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000165'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.92%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.301.32'
$ws.Range('E18').Value = '  -0.83%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '426.56'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.61%  '
$ws.Range('E20').Value = '  -2.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.18'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.72%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.41'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.10'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.76%  '
$ws.Range('E24').Value = '  +0.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.68'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Value = '3.472.08'
$ws.Range('E26').Value = '  -0.51%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.516'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.21%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.206'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.45%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000115'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.58%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.04'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.43%  '
$ws.Range('E31').Value = '  -0.22%  '
$ws.Range('E32').Value = '  -1.11%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '22.46'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.42%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.23'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.65'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.53%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.19'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.64%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '160.92'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.48%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.45'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.82'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.58%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '2.861.83'
$ws.Range('E41').Value = '  +1.83%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '26.52'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.34'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.82%  '
$ws.Range('E44').Value = '  -4.24%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0666'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.12%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.80'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.01'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.84%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.34'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.98%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '23.32'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.46%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '314.47'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.00%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0273'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.51%  '
